# Applies the scheduled-runner profit/price recomputation to Mateus_Profits
# (workbook sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), cell-by-cell, per the diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 45
$ws.Cells.Item(6, 9).Value = 45
$ws.Cells.Item(6, 11).Value = 135
$ws.Cells.Item(6, 13).Value = -23
$ws.Cells.Item(15, 8).Value = 10285.6
$ws.Cells.Item(15, 9).Value = 10285.6
$ws.Cells.Item(15, 11).Value = 30856.8
$ws.Cells.Item(15, 13).Value = -30687.8
$ws.Cells.Item(31, 8).Value = 74.5
$ws.Cells.Item(31, 9).Value = 74.5
$ws.Cells.Item(31, 11).Value = 223.5
$ws.Cells.Item(31, 13).Value = 6.5
$ws.Cells.Item(33, 8).Value = 432.15384
$ws.Cells.Item(33, 9).Value = 158.72728
$ws.Cells.Item(33, 11).Value = 158.72728
$ws.Cells.Item(33, 13).Value = 70.27271999999999
$ws.Cells.Item(98, 8).Value = 2121
$ws.Cells.Item(98, 9).Value = 2228.923
$ws.Cells.Item(98, 11).Value = 2228.923
$ws.Cells.Item(98, 13).Value = -730.9229999999998
$ws.Cells.Item(106, 8).Value = 1500
$ws.Cells.Item(106, 9).Value = 1500
$ws.Cells.Item(106, 11).Value = 1500
$ws.Cells.Item(106, 13).Value = -869
$ws.Cells.Item(122, 8).Value = 2121
$ws.Cells.Item(122, 9).Value = 2228.923
$ws.Cells.Item(122, 11).Value = 6686.768999999999
$ws.Cells.Item(122, 13).Value = -4236.768999999999
$ws.Cells.Item(132, 8).Value = 19516
$ws.Cells.Item(132, 9).Value = 19516
$ws.Cells.Item(132, 11).Value = 58548
$ws.Cells.Item(132, 13).Value = -56018
$ws.Cells.Item(138, 8).Value = 2648.1177
$ws.Cells.Item(138, 10).Value = 2499
$ws.Cells.Item(138, 12).Value = 7497
$ws.Cells.Item(138, 14).Value = -17777

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5687.268
$ws.Cells.Item(61, 9).Value = 4598.3823
$ws.Cells.Item(61, 11).Value = 4598.3823
$ws.Cells.Item(61, 13).Value = -4386.3823
$ws.Cells.Item(92, 8).Value = 44975
$ws.Cells.Item(92, 10).Value = 44975
$ws.Cells.Item(92, 12).Value = 44975
$ws.Cells.Item(92, 14).Value = -49967
$ws.Cells.Item(102, 8).Value = 4322.087
$ws.Cells.Item(102, 9).Value = 3170.45
$ws.Cells.Item(102, 10).Value = 11999.667
$ws.Cells.Item(102, 11).Value = 3170.45
$ws.Cells.Item(102, 12).Value = 11999.667
$ws.Cells.Item(102, 13).Value = -1548.45
$ws.Cells.Item(102, 14).Value = -15243.667
$ws.Cells.Item(132, 8).Value = 4404.9395
$ws.Cells.Item(132, 9).Value = 3294.6428
$ws.Cells.Item(132, 10).Value = 10622.6
$ws.Cells.Item(132, 11).Value = 9883.928400000001
$ws.Cells.Item(132, 12).Value = 31867.8
$ws.Cells.Item(132, 13).Value = -7353.928400000001
$ws.Cells.Item(132, 14).Value = -36927.8
$ws.Cells.Item(136, 8).Value = 5687.268
$ws.Cells.Item(136, 9).Value = 4598.3823
$ws.Cells.Item(136, 11).Value = 13795.1469
$ws.Cells.Item(136, 13).Value = -11245.1469

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 829
$ws.Cells.Item(64, 9).Value = 747.5
$ws.Cells.Item(64, 10).Value = 883.3333
$ws.Cells.Item(64, 11).Value = 747.5
$ws.Cells.Item(64, 12).Value = 883.3333
$ws.Cells.Item(64, 13).Value = -522.5
$ws.Cells.Item(64, 14).Value = -1333.3333
$ws.Cells.Item(67, 8).Value = 829
$ws.Cells.Item(67, 9).Value = 747.5
$ws.Cells.Item(67, 10).Value = 883.3333
$ws.Cells.Item(67, 11).Value = 747.5
$ws.Cells.Item(67, 12).Value = 883.3333
$ws.Cells.Item(67, 13).Value = 32.5
$ws.Cells.Item(67, 14).Value = -2443.3333
$ws.Cells.Item(86, 8).Value = 2124.6316
$ws.Cells.Item(86, 9).Value = 1877.9333
$ws.Cells.Item(86, 10).Value = 3049.75
$ws.Cells.Item(86, 11).Value = 1877.9333
$ws.Cells.Item(86, 12).Value = 3049.75
$ws.Cells.Item(86, 13).Value = -754.9332999999999
$ws.Cells.Item(86, 14).Value = -5295.75
$ws.Cells.Item(89, 8).Value = 2124.6316
$ws.Cells.Item(89, 9).Value = 1877.9333
$ws.Cells.Item(89, 10).Value = 3049.75
$ws.Cells.Item(89, 11).Value = 9389.666499999999
$ws.Cells.Item(89, 12).Value = 15248.75
$ws.Cells.Item(89, 13).Value = -3773.666499999999
$ws.Cells.Item(89, 14).Value = -26480.75
$ws.Cells.Item(105, 8).Value = 2561.75
$ws.Cells.Item(105, 9).Value = 2213.4285
$ws.Cells.Item(105, 10).Value = 5000
$ws.Cells.Item(105, 11).Value = 2213.4285
$ws.Cells.Item(105, 12).Value = 5000
$ws.Cells.Item(105, 13).Value = -466.4285
$ws.Cells.Item(105, 14).Value = -8494
$ws.Cells.Item(134, 8).Value = 3423.973
$ws.Cells.Item(134, 9).Value = 3470.375
$ws.Cells.Item(134, 11).Value = 10411.125
$ws.Cells.Item(134, 13).Value = -7876.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3898
$ws.Cells.Item(16, 9).Value = 3134.8
$ws.Cells.Item(16, 10).Value = 5170
$ws.Cells.Item(16, 11).Value = 3134.8
$ws.Cells.Item(16, 12).Value = 5170
$ws.Cells.Item(16, 13).Value = -2847.8
$ws.Cells.Item(16, 14).Value = -5744
$ws.Cells.Item(105, 8).Value = 1863
$ws.Cells.Item(105, 9).Value = 1863
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 1863
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -116
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 471.96774
$ws.Cells.Item(107, 9).Value = 422.1111
$ws.Cells.Item(107, 11).Value = 422.1111
$ws.Cells.Item(107, 13).Value = 1497.8889
$ws.Cells.Item(113, 8).Value = 3898
$ws.Cells.Item(113, 9).Value = 3134.8
$ws.Cells.Item(113, 10).Value = 5170
$ws.Cells.Item(113, 11).Value = 3134.8
$ws.Cells.Item(113, 12).Value = 5170
$ws.Cells.Item(113, 13).Value = -964.8000000000002
$ws.Cells.Item(113, 14).Value = -9510

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 646.55554
$ws.Cells.Item(113, 9).Value = 480
$ws.Cells.Item(113, 10).Value = 667.375
$ws.Cells.Item(113, 11).Value = 1440
$ws.Cells.Item(113, 12).Value = 2002.125
$ws.Cells.Item(113, 13).Value = 730
$ws.Cells.Item(113, 14).Value = -6342.125
$ws.Cells.Item(132, 8).Value = 1423.25
$ws.Cells.Item(132, 9).Value = 1423.25
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 12809.25
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -10279.25
$ws.Cells.Item(132, 14).ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 874.375
$ws.Cells.Item(14, 9).Value = 248.33333
$ws.Cells.Item(14, 10).Value = 2752.5
$ws.Cells.Item(14, 11).Value = 248.33333
$ws.Cells.Item(14, 12).Value = 2752.5
$ws.Cells.Item(14, 13).Value = -80.33332999999999
$ws.Cells.Item(14, 14).Value = -3088.5
$ws.Cells.Item(97, 8).Value = 3632.389
$ws.Cells.Item(97, 9).Value = 758.93335
$ws.Cells.Item(97, 11).Value = 758.93335
$ws.Cells.Item(97, 13).Value = -262.93335
$ws.Cells.Item(107, 8).Value = 363
$ws.Cells.Item(107, 9).Value = 289.45456
$ws.Cells.Item(107, 11).Value = 289.45456
$ws.Cells.Item(107, 13).Value = 1630.54544
$ws.Cells.Item(122, 8).Value = 2608.3684
$ws.Cells.Item(122, 9).Value = 2381.875
$ws.Cells.Item(122, 11).Value = 7145.625
$ws.Cells.Item(122, 13).Value = -4695.625
$ws.Cells.Item(132, 8).Value = 3095.0312
$ws.Cells.Item(132, 9).Value = 2637.8845
$ws.Cells.Item(132, 10).Value = 5076
$ws.Cells.Item(132, 11).Value = 7913.6535
$ws.Cells.Item(132, 12).Value = 15228
$ws.Cells.Item(132, 13).Value = -5383.6535
$ws.Cells.Item(132, 14).Value = -20288

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 3400
$ws.Cells.Item(19, 9).Value = 3400
$ws.Cells.Item(19, 11).Value = 3400
$ws.Cells.Item(19, 13).Value = -3230
$ws.Cells.Item(61, 8).Value = 53803.844
$ws.Cells.Item(61, 9).Value = 53803.844
$ws.Cells.Item(61, 11).Value = 53803.844
$ws.Cells.Item(61, 13).Value = -53601.844
$ws.Cells.Item(113, 8).Value = 53803.844
$ws.Cells.Item(113, 9).Value = 53803.844
$ws.Cells.Item(113, 11).Value = 53803.844
$ws.Cells.Item(113, 13).Value = -51633.844
$ws.Cells.Item(132, 8).Value = 8570.343000000001
$ws.Cells.Item(132, 10).Value = 7623.5
$ws.Cells.Item(132, 12).Value = 22870.5
$ws.Cells.Item(132, 14).Value = -27930.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 3200
$ws.Cells.Item(23, 9).Value = 2000
$ws.Cells.Item(23, 11).Value = 2000
$ws.Cells.Item(23, 13).Value = -1771
$ws.Cells.Item(100, 8).Value = 1246.1428
$ws.Cells.Item(100, 9).Value = 1241.091
$ws.Cells.Item(100, 11).Value = 2482.182
$ws.Cells.Item(100, 13).Value = -1941.182
$ws.Cells.Item(132, 8).Value = 2270.3215
$ws.Cells.Item(132, 9).Value = 2429.1155
$ws.Cells.Item(132, 10).Value = 206
$ws.Cells.Item(132, 11).Value = 7287.3465
$ws.Cells.Item(132, 12).Value = 618
$ws.Cells.Item(132, 13).Value = -4757.3465
$ws.Cells.Item(132, 14).Value = -5678
$ws.Cells.Item(136, 8).Value = 6098.3125
$ws.Cells.Item(136, 9).Value = 5274.846
$ws.Cells.Item(136, 11).Value = 15824.538
$ws.Cells.Item(136, 13).Value = -13274.538
